$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix casing of the metadata4Ing header labels -> metadata4ing
$ws.Range("B1").Value = "metadata4ing_IRI"
$ws.Range("C1").Value = "metadata4ing_DESC"

# Add the new SBO_DEF column (F)
$ws.Range("F1").Value = "SBO_DEF"
$ws.Range("F2").Value = "[]"

# Match the header formatting used by the other header cells (bold/border style)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
